$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before J (so old J becomes K), shifting existing data right.
$ws.Columns.Item(10).Insert()

# Set header for the newly inserted column J.
$ws.Range("J1").Value = "Personensteuer"

# Update selection to match the target state.
$ws.Range("J7").Select()
